# "varun kumar a paid the fees" — remove her row from the pending list.
# She is row 2 (Reg.No "18bec048", Name "VARUN KUMAR A") on Sheet1;
# deleting the entire row shifts the remaining students up and keeps the
# S.No column, styles and shared-string table consistent automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the whole row so everything below shifts up one place.
$ws.Rows.Item(2).Delete()

# Selection moved from E15 to E11 in the saved file.
$ws.Range("E11").Select() | Out-Null

# The B1 dropdown's source range shrank by one row along with the data.
$dv = $ws.Range("B1").Validation
$dv.Formula1 = "=`$B`$2:`$B`$1208"
